$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the K (Strikeouts) column values for each game row (G2:G9)
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 0
